$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.4694297723977
$ws.Range("C2").Value = 8.683885546606298
$ws.Range("D2").Value = 8.593108863058685
$ws.Range("F2").Value = 37.79784550607309
$ws.Range("G2").Value = 3.690059463900416
$ws.Range("J2").Value = 10.77896947804523
$ws.Range("K2").Value = 11.3506856861423
$ws.Range("L2").Value = 11.44264638805152
$ws.Range("N2").Value = 21.05059790501653
$ws.Range("O2").Value = 28.88443165949083
$ws.Range("B3").Value = 15.28884554334859
$ws.Range("C3").Value = 8.675411024714535
$ws.Range("D3").Value = 8.576835544423885
$ws.Range("F3").Value = 37.87246794870862
$ws.Range("G3").Value = 3.691879701003637
$ws.Range("J3").Value = 10.80086028108293
$ws.Range("K3").Value = 11.22063055540641
$ws.Range("L3").Value = 11.44061519637718
$ws.Range("N3").Value = 21.11144189968054
$ws.Range("O3").Value = 28.96169543763779
$ws.Range("B4").Value = 15.17978464295783
$ws.Range("C4").Value = 8.670320227093759
$ws.Range("D4").Value = 8.568063319501388
$ws.Range("F4").Value = 37.92544024266504
$ws.Range("G4").Value = 3.693057583477803
$ws.Range("J4").Value = 10.8152231293014
$ws.Range("K4").Value = 11.1418983205292
$ws.Range("L4").Value = 11.4408546834273
$ws.Range("N4").Value = 21.15055307974167
$ws.Range("O4").Value = 29.01406356972372
$ws.Range("B5").Value = 15.1358480352021
$ws.Range("C5").Value = 8.668274516789051
$ws.Range("D5").Value = 8.564798025370035
$ws.Range("F5").Value = 37.94882413945309
$ws.Range("G5").Value = 3.693552776935203
$ws.Range("J5").Value = 10.8213083720723
$ws.Range("K5").Value = 11.11013066652833
$ws.Range("L5").Value = 11.44132746088841
$ws.Range("N5").Value = 21.16693322714111
$ws.Range("O5").Value = 29.03664223690106
$ws.Range("B6").Value = 15.1285844158371
$ws.Range("C6").Value = 8.667936584532143
$ws.Range("D6").Value = 8.564274593118473
$ws.Range("F6").Value = 37.95281551411082
$ws.Range("G6").Value = 3.693635922649992
$ws.Range("J6").Value = 10.82233286410175
$ws.Range("K6").Value = 11.1048757689132
$ws.Range("L6").Value = 11.44142866267839
$ws.Range("N6").Value = 21.16967987319547
$ws.Range("O6").Value = 29.04046616503287
$ws.Range("B7").Value = 15.17918998282022
$ws.Range("C7").Value = 8.670292520509047
$ws.Range("D7").Value = 8.568018026107046
$ws.Range("F7").Value = 37.92574833143428
$ws.Range("G7").Value = 3.693064200236455
$ws.Range("J7").Value = 10.81530425595258
$ws.Range("K7").Value = 11.14146856507205
$ws.Range("L7").Value = 11.44085953865709
$ws.Range("N7").Value = 21.15077219669419
$ws.Range("O7").Value = 29.01436306126651
$ws.Range("B8").Value = 15.40681671279406
$ws.Range("C8").Value = 8.680940473442691
$ws.Range("D8").Value = 8.587246400607089
$ws.Range("F8").Value = 37.82208965750536
$ws.Range("G8").Value = 3.690674604923879
$ws.Range("J8").Value = 10.78632638846368
$ws.Range("K8").Value = 11.30563043755872
$ws.Range("L8").Value = 11.44163832577275
$ws.Range("N8").Value = 21.07121391599029
$ws.Range("O8").Value = 28.91004889906328
$ws.Range("B9").Value = 15.8654228435011
$ws.Range("C9").Value = 8.702692359880638
$ws.Range("D9").Value = 8.634501218719317
$ws.Range("F9").Value = 37.67564068144549
$ws.Range("G9").Value = 3.686464545824408
$ws.Range("J9").Value = 10.73679392798464
$ws.Range("K9").Value = 11.63493397923706
$ws.Range("L9").Value = 11.4548980996567
$ws.Range("N9").Value = 20.92904761927961
$ws.Range("O9").Value = 28.74462896259673
$ws.Range("B10").Value = 16.20685124540077
$ws.Range("C10").Value = 8.719175978805309
$ws.Range("D10").Value = 8.674855284577243
$ws.Range("F10").Value = 37.60275842662842
$ws.Range("G10").Value = 3.683658585162745
$ws.Range("J10").Value = 10.70481948861063
$ws.Range("K10").Value = 11.87932627443603
$ws.Range("L10").Value = 11.47170203203484
$ws.Range("N10").Value = 20.83295450368857
$ws.Range("O10").Value = 28.64700680276357
$ws.Range("B11").Value = 16.36252373966465
$ws.Range("C11").Value = 8.72677775661154
$ws.Range("D11").Value = 8.694393801973076
$ws.Range("F11").Value = 37.57714871102992
$ws.Range("G11").Value = 3.682443806053085
$ws.Range("J11").Value = 10.69122649926382
$ws.Range("K11").Value = 11.99060702596748
$ws.Range("L11").Value = 11.48085820715746
$ws.Range("N11").Value = 20.79103559702131
$ws.Range("O11").Value = 28.60779767399776
$ws.Range("B12").Value = 16.42147224510433
$ws.Range("C12").Value = 8.729670644030923
$ws.Range("D12").Value = 8.701958399303136
$ws.Range("F12").Value = 37.5685359290612
$ws.Range("G12").Value = 3.681992620608324
$ws.Range("J12").Value = 10.68621566714397
$ws.Range("K12").Value = 12.03272571783718
$ws.Range("L12").Value = 11.48454074415533
$ws.Range("N12").Value = 20.7754186796307
$ws.Range("O12").Value = 28.59369855763495
$ws.Range("B13").Value = 16.40877750942814
$ws.Range("C13").Value = 8.729046984330779
$ws.Range("D13").Value = 8.700321921432133
$ws.Range("F13").Value = 37.57034258946836
$ws.Range("G13").Value = 3.682089399643888
$ws.Range("J13").Value = 10.68728877339116
$ws.Range("K13").Value = 12.02365620090058
$ws.Range("L13").Value = 11.4837381040335
$ws.Range("N13").Value = 20.77877065652944
$ws.Range("O13").Value = 28.59670175152862
$ws.Range("B14").Value = 16.36737379066634
$ws.Range("C14").Value = 8.727015473857511
$ws.Range("D14").Value = 8.695012847367
$ws.Range("F14").Value = 37.57641838667554
$ws.Range("G14").Value = 3.682406510130205
$ws.Range("J14").Value = 10.69081152098905
$ws.Range("K14").Value = 11.99407277935692
$ws.Range("L14").Value = 11.48115686957404
$ws.Range("N14").Value = 20.78974564252298
$ws.Range("O14").Value = 28.60662272520747
$ws.Range("B15").Value = 16.34201111349037
$ws.Range("C15").Value = 8.725772954224658
$ws.Range("D15").Value = 8.691782352194775
$ws.Range("F15").Value = 37.58028129210601
$ws.Range("G15").Value = 3.682601897484504
$ws.Range("J15").Value = 10.69298707483816
$ws.Range("K15").Value = 11.97594829716257
$ws.Range("L15").Value = 11.47960376092353
$ws.Range("N15").Value = 20.79650155521044
$ws.Range("O15").Value = 28.61279711140024
$ws.Range("B16").Value = 16.19668032892609
$ws.Range("C16").Value = 8.718681212656941
$ws.Range("D16").Value = 8.673601821884619
$ws.Range("F16").Value = 37.60458392225643
$ws.Range("G16").Value = 3.683739211078251
$ws.Range("J16").Value = 10.70572695106473
$ws.Range("K16").Value = 11.87205283809932
$ws.Range("L16").Value = 11.47113390039693
$ws.Range("N16").Value = 20.83573000997952
$ws.Range("O16").Value = 28.64967391495271
$ws.Range("B17").Value = 16.10757798414328
$ws.Range("C17").Value = 8.714356694708476
$ws.Range("D17").Value = 8.662748358859353
$ws.Range("F17").Value = 37.62142541637385
$ws.Range("G17").Value = 3.684452680020182
$ws.Range("J17").Value = 10.71378607908515
$ws.Range("K17").Value = 11.80831747082956
$ws.Range("L17").Value = 11.46632369697995
$ws.Range("N17").Value = 20.86025413933389
$ws.Range("O17").Value = 28.67362901886735
$ws.Range("B18").Value = 16.05636516473822
$ws.Range("C18").Value = 8.711879105183877
$ws.Range("D18").Value = 8.656617174692501
$ws.Range("F18").Value = 37.63182236617357
$ws.Range("G18").Value = 3.684868855344298
$ws.Range("J18").Value = 10.71851113959513
$ws.Range("K18").Value = 11.77167062178887
$ws.Range("L18").Value = 11.46369944581931
$ws.Range("N18").Value = 20.87452872971526
$ws.Range("O18").Value = 28.68789673664556
$ws.Range("B19").Value = 16.03903326223381
$ws.Range("C19").Value = 8.711041930036004
$ws.Range("D19").Value = 8.654560522377409
$ws.Range("F19").Value = 37.63546455639472
$ws.Range("G19").Value = 3.685010763793515
$ws.Range("J19").Value = 10.7201263761324
$ws.Range("K19").Value = 11.75926580811893
$ws.Range("L19").Value = 11.462835448221
$ws.Range("N19").Value = 20.87939091529761
$ws.Range("O19").Value = 28.69281156481509
$ws.Range("B20").Value = 16.11705968626928
$ws.Range("C20").Value = 8.714816041314929
$ws.Range("D20").Value = 8.663892223799786
$ws.Range("F20").Value = 37.61955910867231
$ws.Range("G20").Value = 3.684376129354978
$ws.Range("J20").Value = 10.71291889430886
$ws.Range("K20").Value = 11.81510122605265
$ws.Range("L20").Value = 11.46682102325125
$ws.Range("N20").Value = 20.85762602562113
$ws.Range("O20").Value = 28.67102830234107
$ws.Range("B21").Value = 16.37953551352286
$ws.Range("C21").Value = 8.727611795513907
$ws.Range("D21").Value = 8.696567785149142
$ws.Range("F21").Value = 37.57460433040842
$ws.Range("G21").Value = 3.682313127921637
$ws.Range("J21").Value = 10.68977310225897
$ws.Range("K21").Value = 12.00276299896026
$ws.Range("L21").Value = 11.48190921535025
$ws.Range("N21").Value = 20.7865150611088
$ws.Range("O21").Value = 28.60368837381459
$ws.Range("B22").Value = 16.55104762228733
$ws.Range("C22").Value = 8.736057392989711
$ws.Range("D22").Value = 8.718887457949219
$ws.Range("F22").Value = 37.55154813417754
$ws.Range("G22").Value = 3.68101625683606
$ws.Range("J22").Value = 10.67544165009943
$ws.Range("K22").Value = 12.12527271782879
$ws.Range("L22").Value = 11.4930240294642
$ws.Range("N22").Value = 20.74153659729675
$ws.Range("O22").Value = 28.56404115158323
$ws.Range("B23").Value = 16.4595282339971
$ws.Range("C23").Value = 8.731542437785587
$ws.Range("D23").Value = 8.706888187652087
$ws.Range("F23").Value = 37.56327504049413
$ws.Range("G23").Value = 3.681703730400181
$ws.Range("J23").Value = 10.68301794660638
$ws.Range("K23").Value = 12.05991133856565
$ws.Range("L23").Value = 11.48697785739445
$ws.Range("N23").Value = 20.76540588576901
$ws.Range("O23").Value = 28.58480213133234
$ws.Range("B24").Value = 16.11277296328374
$ws.Range("C24").Value = 8.714608343783015
$ws.Range("D24").Value = 8.663374743837329
$ws.Range("F24").Value = 37.62040064061789
$ws.Range("G24").Value = 3.684410719241006
$ws.Range("J24").Value = 10.71331066269087
$ws.Range("K24").Value = 11.81203430168351
$ws.Range("L24").Value = 11.46659574207019
$ws.Range("N24").Value = 20.85881364955139
$ws.Range("O24").Value = 28.67220254254229
$ws.Range("B25").Value = 15.74035457517498
$ws.Range("C25").Value = 8.696717551040159
$ws.Range("D25").Value = 8.620713630533816
$ws.Range("F25").Value = 37.70916670088197
$ws.Range("G25").Value = 3.68755283511945
$ws.Range("J25").Value = 10.74941600977243
$ws.Range("K25").Value = 11.54526988411826
$ws.Range("L25").Value = 11.45006375653252
$ws.Range("N25").Value = 20.96603377198828
$ws.Range("O25").Value = 28.78518363954199
